$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the old row 5 ("Funkmast"), shifting the
# existing rows 5-14 down to 6-15.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new lookup value.
$ws.Range("A5").Value = "Freileitungsmast"
$ws.Range("B5").Value = 1251
$ws.Range("C5").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/industrial"

# The row insert landed inside the old multi-cell hyperlink (C3:C14),
# which no longer covers the full C3:C15 data range correctly, so it
# gets dropped (matching the source edit, which leaves only the C2
# hyperlink in place and plain, unlinked URL text in C3:C15).
$ws.Hyperlinks.Item(2).Delete()

# Match the author's final selection.
$ws.Range("C5").Select()
